$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 3 (Away/Road splits), Week 17 data logged
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 559
$wsOff.Range("C3").Value = 385
$wsOff.Range("D3").Value = 117
$wsOff.Range("E3").Value = 63
$wsOff.Range("F3").Value = 15
$wsOff.Range("G3").Value = 4

# Sheet "DEF" - row 3 (Away/Road splits), Week 17 data logged
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 452
$wsDef.Range("C3").Value = 316
$wsDef.Range("D3").Value = 120
$wsDef.Range("E3").Value = 69
